$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for rows 2-121.
# Update all of them to the new date serial 45175 (2023-09-06),
# matching the original numeric storage (date-formatted cell, style preserved).
$ws.Range("C2:C121").Value = 45175
